$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (BGR values as used by Excel's Interior.Color / COM)
$winColor = 32768      # RGB 00008000 (green)   -> "Win"
$neutralColor = 42495  # RGB 00FFA500 (orange)  -> "Neutral"

# --- Row 2: LoadPlist / any,plistlib.readPlistFromString -> reordered set; status Neutral -> Win
$ws.Range("E2").Value2 = "{'plistlib.readPlistFromString', 'any'}"
$ws.Range("F2").Value2 = "Win"
$ws.Range("F2").Interior.Color = $winColor

# --- Row 3: LoadPlist / any -> plistlib.readPlistFromString; status Neutral -> Win
$ws.Range("E3").Value2 = "plistlib.readPlistFromString"
$ws.Range("F3").Value2 = "Win"
$ws.Range("F3").Interior.Color = $winColor

# --- Row 4: GetPlistValue; status Neutral -> Win (E4 unchanged)
$ws.Range("F4").Value2 = "Win"
$ws.Range("F4").Interior.Color = $winColor

# --- Row 5: GetPlistValue; stays Neutral (re-assert orange to keep it consistent with fill renumbering)
$ws.Range("F5").Value2 = "Neutral"
$ws.Range("F5").Interior.Color = $neutralColor

# --- Row 6: GetProgram; set reordered, stays Neutral
$ws.Range("E6").Value2 = "{'Tuple[str]', 'any', 'Tuple[any]'}"
$ws.Range("F6").Value2 = "Neutral"
$ws.Range("F6").Interior.Color = $neutralColor

# --- Row 7: GetProgram; any -> Tuple[str]; status Loss -> Neutral
$ws.Range("E7").Value2 = "Tuple[str]"
$ws.Range("F7").Value2 = "Neutral"
$ws.Range("F7").Interior.Color = $neutralColor

# --- Row 8: HashFile; stays Neutral
$ws.Range("F8").Value2 = "Neutral"
$ws.Range("F8").Interior.Color = $neutralColor

# --- Row 9: HashFile; stays Neutral
$ws.Range("F9").Value2 = "Neutral"
$ws.Range("F9").Interior.Color = $neutralColor

# --- Row 10: GetComment; stays Neutral
$ws.Range("F10").Value2 = "Neutral"
$ws.Range("F10").Interior.Color = $neutralColor

# --- Row 11: GetComment; stays Neutral
$ws.Range("F11").Value2 = "Neutral"
$ws.Range("F11").Interior.Color = $neutralColor

# --- Row 12: LoadPlist; stays Win
$ws.Range("F12").Value2 = "Win"
$ws.Range("F12").Interior.Color = $winColor

# --- Row 13: GetPlistValue; stays Neutral
$ws.Range("F13").Value2 = "Neutral"
$ws.Range("F13").Interior.Color = $neutralColor

# --- Row 14: GetPlistValue; stays Neutral
$ws.Range("F14").Value2 = "Neutral"
$ws.Range("F14").Interior.Color = $neutralColor

# --- Row 15: GetProgram; stays Neutral
$ws.Range("F15").Value2 = "Neutral"
$ws.Range("F15").Interior.Color = $neutralColor

# --- Row 16: HashFile; stays Win
$ws.Range("F16").Value2 = "Win"
$ws.Range("F16").Interior.Color = $winColor

# --- Row 17: GetComment; stays Neutral
$ws.Range("F17").Value2 = "Neutral"
$ws.Range("F17").Interior.Color = $neutralColor

# --- Row 18: GetComment; stays Neutral
$ws.Range("F18").Value2 = "Neutral"
$ws.Range("F18").Interior.Color = $neutralColor

# --- Row 19: summary counts - PyType Wins 1 -> 0 ; Scalpel Wins 2 -> 5
$ws.Range("D19").Value2 = 0
$ws.Range("F19").Value2 = 5

# --- Insert a new row at 20 (pushes old row 20 "Accuracy over PyType" content down to row 21)
$ws.Rows("20").Insert()

# --- New row 20: Scalpel Accuracy summary
$ws.Range("C20").Value2 = "Scalpel Accuracy:"
$ws.Range("D20").Value2 = 100

# --- Row 21 (shifted down from old row 20): update F21 value 200 -> 100
$ws.Range("F21").Value2 = 100
